# Applies the "account_creation" sheet addition described by the diff.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the existing "credentials" sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "account_creation"

# --- Apply the text ("@") number format to the whole data area up front ---
$ws.Range("A1:I3").NumberFormat = "@"

# --- Header row ---
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "firstName"
$ws.Range("C1").Value = "lastName"
$ws.Range("D1").Value = "password"
$ws.Range("E1").Value = "birthDay"
$ws.Range("F1").Value = "birthMonth"
$ws.Range("G1").Value = "birthYear"
$ws.Range("H1").Value = "newsletters"
$ws.Range("I1").Value = "offers"

# --- Name columns for both users (A, B, C) ---
$ws.Range("A2").Value = "Mr"
$ws.Range("A3").Value = "Mrs"
$ws.Range("B2").Value = "John"
$ws.Range("B3").Value = "Jane"
$ws.Range("C2").Value = "Doe"
$ws.Range("C3").Value = "Doe"

# --- newsletters/offers for user 1 (row 2) ---
$ws.Range("H2").Value = "'true"
$ws.Range("I2").Value = "'true"

# --- password / birth date columns for both users (D, E, F, G) ---
$ws.Range("D2").Value = 123456
$ws.Range("D3").Value = "'123456"
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = "'20"
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = "'12"
$ws.Range("G2").Value = 1990
$ws.Range("G3").Value = "'1996"

# --- newsletters/offers for user 2 (row 3) ---
$ws.Range("H3").Value = "'false"
$ws.Range("I3").Value = "'false"

# --- Selection on the new sheet, matching the target workbook state ---
$ws.Range("D4").Select()
